$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.676558375358582
$ws.Range("B1").Value = 1.808444857597351
$ws.Range("C1").Value = 2.058799266815186
$ws.Range("D1").Value = 3.283451557159424
$ws.Range("E1").Value = 3.321225166320801
